$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 346-347, shifting the existing rows 346:349 down to 348:351
$ws.Rows("346:347").Insert()

# New row 346
$ws.Range("A346").Value = 10
$ws.Range("B346").Value = "Vega Modelo de Temuco"
$ws.Range("C346").Value = "La Araucanía"
$ws.Range("D346").Value = 44595
$ws.Range("E346").Value = 9
$ws.Range("F346").Value = 100112032
$ws.Range("G346").Value = "Zapallo italiano"
$ws.Range("H346").Value = "Sin especificar"
$ws.Range("I346").Value = "Primera"
$ws.Range("J346").Value = 20
$ws.Range("K346").Value = 10000
$ws.Range("L346").Value = 10000
$ws.Range("M346").Value = 10000
$ws.Range("N346").Value = "$/caja 36 unidades"
$ws.Range("O346").Value = "Región de La Araucanía"
$ws.Range("P346").Value = 278
$ws.Range("Q346").Value = 36
$ws.Range("R346").Value = "Hortaliza"

# New row 347
$ws.Range("A347").Value = 10
$ws.Range("B347").Value = "Vega Modelo de Temuco"
$ws.Range("C347").Value = "La Araucanía"
$ws.Range("D347").Value = 44595
$ws.Range("E347").Value = 9
$ws.Range("F347").Value = 100112032
$ws.Range("G347").Value = "Zapallo italiano"
$ws.Range("H347").Value = "Sin especificar"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 300
$ws.Range("K347").Value = 12000
$ws.Range("L347").Value = 12000
$ws.Range("M347").Value = 12000
$ws.Range("N347").Value = "$/caja 60 unidades"
$ws.Range("O347").Value = "Región de Arica y Parinacota"
$ws.Range("P347").Value = 200
$ws.Range("Q347").Value = 60
$ws.Range("R347").Value = "Hortaliza"
